# Bump the "code SiteList index" counters (column F) by 2 on every data
# sheet, double the "Conductance" sheet's E1 value, and leave the cell
# selection on each sheet's first cell (mirroring what Excel records after
# a user clicks through each tab from A1), ending on the ESTM Coefficients
# tab which stays the active tab (as in the source workbook).

$wb = $excel.ActiveWorkbook

# Paved
$ws = $wb.Worksheets.Item("Paved")
$ws.Range("F1").Value = 36
$ws.Range("F1").Select()

# Building
$ws = $wb.Worksheets.Item("Building")
$ws.Range("F1").Value = 37
$ws.Range("F1").Select()

# Evergreen
$ws = $wb.Worksheets.Item("Evergreen")
$ws.Range("F1").Value = 38
$ws.Range("F1").Select()

# Decidious
$ws = $wb.Worksheets.Item("Decidious")
$ws.Range("F1").Value = 39
$ws.Range("F1").Select()

# Grass
$ws = $wb.Worksheets.Item("Grass")
$ws.Range("F1").Value = 40
$ws.Range("F1").Select()

# Bare Soil
$ws = $wb.Worksheets.Item("Bare Soil")
$ws.Range("F1").Value = 41
$ws.Range("F1").Select()

# Water
$ws = $wb.Worksheets.Item("Water")
$ws.Range("F1").Value = 42
$ws.Range("F1").Select()

# Conductance
$ws = $wb.Worksheets.Item("Conductance")
$ws.Range("E1").Value = 200
$ws.Range("F1").Value = 47
$ws.Range("F1").Select()

# Snow
$ws = $wb.Worksheets.Item("Snow")
$ws.Range("F1").Value = 48
$ws.Range("F1").Select()

# Snow clearing
$ws = $wb.Worksheets.Item("Snow clearing")
$ws.Range("F1").Value = 49
$ws.Range("F2").Value = 50
$ws.Range("F2").Select()

# Anthropogenic
$ws = $wb.Worksheets.Item("Anthropogenic")
$ws.Range("F1").Value = 51
$ws.Range("F1").Select()

# Energy
$ws = $wb.Worksheets.Item("Energy")
$ws.Range("F1").Value = 52
$ws.Range("F2").Value = 53
$ws.Range("F2").Select()

# Irrigation
$ws = $wb.Worksheets.Item("Irrigation")
$ws.Range("F1").Value = 56
$ws.Range("F1").Select()

# Water Use (Manual)
$ws = $wb.Worksheets.Item("Water Use (Manual)")
$ws.Range("F1").Value = 57
$ws.Range("F2").Value = 58
$ws.Range("F2").Select()

# Water Use (Automatic)
$ws = $wb.Worksheets.Item("Water Use (Automatic)")
$ws.Range("F1").Value = 59
$ws.Range("F2").Value = 60
$ws.Range("F2").Select()

# ESTM Coefficients (stays the active/visible tab, like in the source file)
$ws = $wb.Worksheets.Item("ESTM Coefficients")
$ws.Range("F1").Value = 87
$ws.Range("F5").Select()
